# Scheduled-runner style refresh of market-price / leve-profit figures
# across several crafting-job sheets (currentAveragePrice*, LevePrice*,
# LeveProfit* columns H:N). Values below are the refreshed raw data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 14159.167
$ws.Range("I127").Value = 1989.1666
$ws.Range("J127").Value = 38499.168
$ws.Range("K127").Value = 5967.4998
$ws.Range("L127").Value = 115497.504
$ws.Range("M127").Value = -1007.4998
$ws.Range("N127").Value = -125417.504

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1852.9524
$ws.Range("I2").Value = 1853.3549
$ws.Range("J2").Value = 1851.8182
$ws.Range("K2").Value = 1853.3549
$ws.Range("L2").Value = 1851.8182
$ws.Range("M2").Value = -1740.3549
$ws.Range("N2").Value = -2077.8182

$ws.Range("H61").Value = 3818.55
$ws.Range("I61").Value = 2995.4517
$ws.Range("J61").Value = 6653.6665
$ws.Range("K61").Value = 2995.4517
$ws.Range("L61").Value = 6653.6665
$ws.Range("M61").Value = -2783.4517
$ws.Range("N61").Value = -7077.6665

$ws.Range("H74").Value = 1975.9333
$ws.Range("I74").Value = 1575.5385
$ws.Range("J74").Value = 4578.5
$ws.Range("K74").Value = 1575.5385
$ws.Range("L74").Value = 4578.5
$ws.Range("M74").Value = -701.5385000000001
$ws.Range("N74").Value = -6326.5

$ws.Range("H77").Value = 1975.9333
$ws.Range("I77").Value = 1575.5385
$ws.Range("J77").Value = 4578.5
$ws.Range("K77").Value = 7877.692500000001
$ws.Range("L77").Value = 22892.5
$ws.Range("M77").Value = -3509.692500000001
$ws.Range("N77").Value = -31628.5

$ws.Range("H116").Value = 1852.9524
$ws.Range("I116").Value = 1853.3549
$ws.Range("J116").Value = 1851.8182
$ws.Range("K116").Value = 1853.3549
$ws.Range("L116").Value = 1851.8182
$ws.Range("M116").Value = 440.6451
$ws.Range("N116").Value = -6439.8182

$ws.Range("H136").Value = 3818.55
$ws.Range("I136").Value = 2995.4517
$ws.Range("J136").Value = 6653.6665
$ws.Range("K136").Value = 8986.355100000001
$ws.Range("L136").Value = 19960.9995
$ws.Range("M136").Value = -6436.355100000001
$ws.Range("N136").Value = -25060.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1852.9524
$ws.Range("I3").Value = 1853.3549
$ws.Range("J3").Value = 1851.8182
$ws.Range("K3").Value = 1853.3549
$ws.Range("L3").Value = 1851.8182
$ws.Range("M3").Value = -1739.3549
$ws.Range("N3").Value = -2079.8182

$ws.Range("H86").Value = 5027.722
$ws.Range("I86").Value = 3456.2856
$ws.Range("J86").Value = 6027.727
$ws.Range("K86").Value = 3456.2856
$ws.Range("L86").Value = 6027.727
$ws.Range("M86").Value = -2333.2856
$ws.Range("N86").Value = -8273.726999999999

$ws.Range("H89").Value = 5027.722
$ws.Range("I89").Value = 3456.2856
$ws.Range("J89").Value = 6027.727
$ws.Range("K89").Value = 17281.428
$ws.Range("L89").Value = 30138.635
$ws.Range("M89").Value = -11665.428
$ws.Range("N89").Value = -41370.63499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 6813.0586
$ws.Range("I94").Value = 16978.666
$ws.Range("J94").Value = 1268.1818
$ws.Range("K94").Value = 16978.666
$ws.Range("L94").Value = 1268.1818
$ws.Range("M94").Value = -16527.666
$ws.Range("N94").Value = -2170.1818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2584.4546
$ws.Range("I3").Value = 2718.4285
$ws.Range("J3").Value = 2350
$ws.Range("K3").Value = 8155.2855
$ws.Range("L3").Value = 7050
$ws.Range("M3").Value = -8043.2855
$ws.Range("N3").Value = -7274

$ws.Range("H12").Value = 175.82353
$ws.Range("J12").Value = 201.1
$ws.Range("L12").Value = 603.3
$ws.Range("N12").Value = -949.3

$ws.Range("H63").Value = 4690.2856
$ws.Range("I63").Value = 3732.8
$ws.Range("K63").Value = 11198.4
$ws.Range("M63").Value = -10449.4

$ws.Range("H66").Value = 4690.2856
$ws.Range("I66").Value = 3732.8
$ws.Range("K66").Value = 33595.2
$ws.Range("M66").Value = -29851.2

$ws.Range("H69").Value = 6000
$ws.Range("J69").Value = 6000
$ws.Range("L69").Value = 18000
$ws.Range("N69").Value = -19622

$ws.Range("H72").Value = 6000
$ws.Range("J72").Value = 6000
$ws.Range("L72").Value = 54000
$ws.Range("N72").Value = -62112

$ws.Range("H75").Value = 1715.0625
$ws.Range("I75").Value = 628
$ws.Range("J75").Value = 1965.9231
$ws.Range("K75").Value = 1884
$ws.Range("L75").Value = 5897.7693
$ws.Range("M75").Value = -886
$ws.Range("N75").Value = -7893.7693

$ws.Range("H78").Value = 1715.0625
$ws.Range("I78").Value = 628
$ws.Range("J78").Value = 1965.9231
$ws.Range("K78").Value = 5652
$ws.Range("L78").Value = 17693.3079
$ws.Range("M78").Value = -660
$ws.Range("N78").Value = -27677.3079

$ws.Range("H88").Value = 4250
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 4250
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 12750
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -13606

$ws.Range("H91").Value = 4250
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 4250
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 12750
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -15714

$ws.Range("H123").Value = 2015
$ws.Range("I123").Value = 2015
$ws.Range("K123").Value = 6045
$ws.Range("M123").Value = -3595

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 47117.76
$ws.Range("I80").Value = 71834.625
$ws.Range("J80").Value = 3176.6667
$ws.Range("K80").Value = 71834.625
$ws.Range("L80").Value = 3176.6667
$ws.Range("M80").Value = -70836.625
$ws.Range("N80").Value = -5172.6667

$ws.Range("H83").Value = 47117.76
$ws.Range("I83").Value = 71834.625
$ws.Range("J83").Value = 3176.6667
$ws.Range("K83").Value = 359173.125
$ws.Range("L83").Value = 15883.3335
$ws.Range("M83").Value = -354181.125
$ws.Range("N83").Value = -25867.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 6093.8335
$ws.Range("I82").Value = 3285.4285
$ws.Range("J82").Value = 10025.6
$ws.Range("K82").Value = 3285.4285
$ws.Range("L82").Value = 10025.6
$ws.Range("M82").Value = -2924.4285
$ws.Range("N82").Value = -10747.6

$ws.Range("H85").Value = 6093.8335
$ws.Range("I85").Value = 3285.4285
$ws.Range("J85").Value = 10025.6
$ws.Range("K85").Value = 3285.4285
$ws.Range("L85").Value = 10025.6
$ws.Range("M85").Value = -2037.4285
$ws.Range("N85").Value = -12521.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2991.8333
$ws.Range("I96").Value = 3009.4
$ws.Range("J96").Value = 2904
$ws.Range("K96").Value = 3009.4
$ws.Range("L96").Value = 2904
$ws.Range("M96").Value = -1636.4
$ws.Range("N96").Value = -5650

$ws.Range("H132").Value = 1275.0476
$ws.Range("I132").Value = 726.6177
$ws.Range("J132").Value = 3605.875
$ws.Range("K132").Value = 2179.8531
$ws.Range("L132").Value = 10817.625
$ws.Range("M132").Value = 350.1468999999997
$ws.Range("N132").Value = -15877.625

$ws.Range("H136").Value = 2198.6272
$ws.Range("I136").Value = 2165.111
$ws.Range("J136").Value = 2251.087
$ws.Range("K136").Value = 6495.333
$ws.Range("L136").Value = 6753.261
$ws.Range("M136").Value = -3945.333
$ws.Range("N136").Value = -11853.261
